$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-8, 0, 6.5467, 30, 0, 250, 112, 0.16, 0.158),
    @(-7, 0, 7.1989, 30, 0, 250, 112, 0.16, 0.158),
    @(-6, 0, 7.5776, 30, 0, 250, 112, 0.16, 0.158),
    @(-5, 0, 7.544, 30, 0, 250, 112, 0.16, 0.158),
    @(-4, 0, 7.3532, 30, 0, 250, 112, 0.16, 0.158),
    @(-3, 0, 7.0982, 30, 0, 250, 112, 0.16, 0.158),
    @(-2, 0, 5.9479, 30, 0, 250, 112, 0.16, 0.158),
    @(-1, 0, 4.4079, 30, 0, 250, 112, 0.16, 0.158),
    @(0, 0, 2.3276, 30, 0, 250, 112, 0.16, 0.158),
    @(1, 0, 0.6317, 30, 0, 250, 112, 0.16, 0.158),
    @(2, 0, -2.1316, 30, 0, 250, 112, 0.16, 0.158),
    @(3, 0, -5.237, 30, 0, 250, 112, 0.16, 0.158),
    @(-8, 3, 6.4715, 30, 0, 250, 112, 0.16, 0.158),
    @(-3, 3, 7.1359, 30, 0, 250, 112, 0.16, 0.158),
    @(-2, 3, 6.2015, 30, 0, 250, 112, 0.16, 0.158),
    @(-1, 3, 5.3492, 30, 0, 250, 112, 0.16, 0.158),
    @(0, 3, 3.0631, 30, 0, 250, 112, 0.16, 0.158),
    @(1, 3, 0.40723, 30, 0, 250, 112, 0.16, 0.158),
    @(2, 3, -0.75966, 30, 0, 250, 112, 0.16, 0.158),
    @(3, 3, -3.4912, 30, 0, 250, 112, 0.16, 0.158),
    @(-8, 6, 6.3803, 30, 0, 250, 112, 0.16, 0.158),
    @(3, 6, -2.2279, 30, 0, 250, 112, 0.16, 0.158),
    @(-8, 9, 6.2964, 30, 0, 250, 112, 0.16, 0.158),
    @(-6, 9, 7.4468, 30, 0, 250, 112, 0.16, 0.158),
    @(0, 9, 4.5422, 30, 0, 250, 112, 0.16, 0.158),
    @(1, 9, 1.9925, 30, 0, 250, 112, 0.16, 0.158),
    @(2, 9, 0.34576, 30, 0, 250, 112, 0.16, 0.158),
    @(3, 9, -1.0627, 30, 0, 250, 112, 0.16, 0.158),
    @(-8, 12, 6.1934, 30, 0, 250, 112, 0.16, 0.158),
    @(-3, 30, 7.3255, 30, 0, 250, 112, 0.16, 0.158),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
